$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update summary counts: Cant. Trabajadores (workers) and Cant. Periodos (periods)
$ws.Range("C13").Value2 = 28
$ws.Range("F13").Value2 = 2

# Update the detail table (rows 16-44): replace period 2507 batch with
# the combined 2507 (LUISA carry-over) + new 2508 batch data.
$ws.Range("C16").Value2 = '45547986'
$ws.Range("D16").Value2 = 'LUISA CARMELA VASQUEZ BELTRAN'
$ws.Range("E16").Value2 = '2507'
$ws.Range("F16").Value2 = 56940
$ws.Range("G16").Value2 = 1423500
$ws.Range("C17").Value2 = '1045670210'
$ws.Range("D17").Value2 = 'DIEGO ALEXANDER DIAZ ORTIZ'
$ws.Range("E17").Value2 = '2508'
$ws.Range("F17").Value2 = 280000
$ws.Range("G17").Value2 = 7000000
$ws.Range("C18").Value2 = '1128063151'
$ws.Range("D18").Value2 = 'KATYA MILENIS RODRIGUEZ DIAZ'
$ws.Range("E18").Value2 = '2508'
$ws.Range("F18").Value2 = 60940
$ws.Range("G18").Value2 = 1523500
$ws.Range("C19").Value2 = '45499877'
$ws.Range("D19").Value2 = 'MARLY PEREZ LOMBANA'
$ws.Range("E19").Value2 = '2508'
$ws.Range("F19").Value2 = 56940
$ws.Range("G19").Value2 = 1423500
$ws.Range("C20").Value2 = '1128062374'
$ws.Range("D20").Value2 = 'VICTOR ALFONSO MARIMON AGUASLIMPIAS'
$ws.Range("E20").Value2 = '2508'
$ws.Range("F20").Value2 = 72000
$ws.Range("G20").Value2 = 1800000
$ws.Range("C21").Value2 = '45541075'
$ws.Range("D21").Value2 = 'YANIRIS MERCEDES BOLIVAR VELEZ'
$ws.Range("E21").Value2 = '2508'
$ws.Range("F21").Value2 = 92400
$ws.Range("G21").Value2 = 2310000
$ws.Range("C22").Value2 = '15621041'
$ws.Range("D22").Value2 = 'DUBAN JOSE OLEA FERNANDEZ'
$ws.Range("E22").Value2 = '2508'
$ws.Range("F22").Value2 = 56940
$ws.Range("G22").Value2 = 1423500
$ws.Range("C23").Value2 = '33100838'
$ws.Range("D23").Value2 = 'ANA ISABEL MORENO CARMONA'
$ws.Range("E23").Value2 = '2508'
$ws.Range("F23").Value2 = 56940
$ws.Range("G23").Value2 = 1423500
$ws.Range("C24").Value2 = '73168948'
$ws.Range("D24").Value2 = 'GERONIMO DE LA ROSA RIVAS'
$ws.Range("E24").Value2 = '2508'
$ws.Range("F24").Value2 = 72000
$ws.Range("G24").Value2 = 1800000
$ws.Range("C25").Value2 = '1003504730'
$ws.Range("D25").Value2 = 'NEIDER ENRIQUE CONTRERAS MENDEZ'
$ws.Range("E25").Value2 = '2508'
$ws.Range("F25").Value2 = 88000
$ws.Range("G25").Value2 = 2200000
$ws.Range("C26").Value2 = '45547986'
$ws.Range("D26").Value2 = 'LUISA CARMELA VASQUEZ BELTRAN'
$ws.Range("E26").Value2 = '2508'
$ws.Range("F26").Value2 = 56940
$ws.Range("G26").Value2 = 1423500
$ws.Range("C27").Value2 = '1143394628'
$ws.Range("D27").Value2 = 'YAMILE RAMOS GODIN'
$ws.Range("E27").Value2 = '2508'
$ws.Range("F27").Value2 = 56940
$ws.Range("G27").Value2 = 1423500
$ws.Range("C28").Value2 = '1143326648'
$ws.Range("D28").Value2 = 'JORGY FLOREZ LARA'
$ws.Range("E28").Value2 = '2508'
$ws.Range("F28").Value2 = 56940
$ws.Range("G28").Value2 = 1423500
$ws.Range("C29").Value2 = '1047400272'
$ws.Range("D29").Value2 = 'STEFANI CASTRO MARTINEZ'
$ws.Range("E29").Value2 = '2508'
$ws.Range("F29").Value2 = 92400
$ws.Range("G29").Value2 = 2310000
$ws.Range("C30").Value2 = '1047502383'
$ws.Range("D30").Value2 = 'JAIME SILVA GONZALEZ'
$ws.Range("E30").Value2 = '2508'
$ws.Range("F30").Value2 = 71400
$ws.Range("G30").Value2 = 1785000
$ws.Range("C31").Value2 = '1149445719'
$ws.Range("D31").Value2 = 'JOSE ANTONIO ORTEGA RODRIGUEZ'
$ws.Range("E31").Value2 = '2508'
$ws.Range("F31").Value2 = 72000
$ws.Range("G31").Value2 = 1800000
$ws.Range("C32").Value2 = '1143414865'
$ws.Range("D32").Value2 = 'SERGIO ANDRES VARGAS SOTOMAYOR'
$ws.Range("E32").Value2 = '2508'
$ws.Range("F32").Value2 = 56940
$ws.Range("G32").Value2 = 1423500
$ws.Range("C33").Value2 = '1042584373'
$ws.Range("D33").Value2 = 'YENIFER PADILLA CUESTA'
$ws.Range("E33").Value2 = '2508'
$ws.Range("F33").Value2 = 56940
$ws.Range("G33").Value2 = 1423500
$ws.Range("C34").Value2 = '1003202611'
$ws.Range("D34").Value2 = 'ALFRED ANDRES QUINTERO QUINTERO'
$ws.Range("E34").Value2 = '2508'
$ws.Range("F34").Value2 = 88000
$ws.Range("G34").Value2 = 2200000
$ws.Range("C35").Value2 = '1047447759'
$ws.Range("D35").Value2 = 'KELLY DEL CARMEN CONTRERAS MEDINA'
$ws.Range("E35").Value2 = '2508'
$ws.Range("F35").Value2 = 56940
$ws.Range("G35").Value2 = 1423500
$ws.Range("C36").Value2 = '1007523433'
$ws.Range("D36").Value2 = 'CLAUDIA PATRICIA SALGADO MEDRANO'
$ws.Range("E36").Value2 = '2508'
$ws.Range("F36").Value2 = 56940
$ws.Range("G36").Value2 = 1423500
$ws.Range("C37").Value2 = '1032452616'
$ws.Range("D37").Value2 = 'ALVARO ANDRES CAMARGO TOVAR'
$ws.Range("E37").Value2 = '2508'
$ws.Range("F37").Value2 = 440000
$ws.Range("G37").Value2 = 11000000
$ws.Range("C38").Value2 = '1057607305'
$ws.Range("D38").Value2 = 'GERALDINE NATALIA SALAMANCA MUÑOZ'
$ws.Range("E38").Value2 = '2508'
$ws.Range("F38").Value2 = 71400
$ws.Range("G38").Value2 = 1785000
$ws.Range("C39").Value2 = '1033807078'
$ws.Range("D39").Value2 = 'LIZETH NATALIA BERNAL GAMBA'
$ws.Range("E39").Value2 = '2508'
$ws.Range("F39").Value2 = 104000
$ws.Range("G39").Value2 = 2600000
$ws.Range("C40").Value2 = '1041973482'
$ws.Range("D40").Value2 = 'ANDRES DAVID QUINTERO OSPINO'
$ws.Range("E40").Value2 = '2508'
$ws.Range("F40").Value2 = 56940
$ws.Range("G40").Value2 = 1423500
$ws.Range("C41").Value2 = '1002201240'
$ws.Range("D41").Value2 = 'PAULA ANDREA PATERNINA MUÑOZ'
$ws.Range("E41").Value2 = '2508'
$ws.Range("F41").Value2 = 56940
$ws.Range("G41").Value2 = 1423500
$ws.Range("C42").Value2 = '1052096113'
$ws.Range("D42").Value2 = 'MARIO RAFAEL VERGARA CASTRO'
$ws.Range("E42").Value2 = '2508'
$ws.Range("F42").Value2 = 56940
$ws.Range("G42").Value2 = 1423500
$ws.Range("C43").Value2 = '1007314639'
$ws.Range("D43").Value2 = 'JHONATANN GUERRERO AGAMEZ'
$ws.Range("E43").Value2 = '2508'
$ws.Range("F43").Value2 = 56940
$ws.Range("G43").Value2 = 1423500
$ws.Range("C44").Value2 = '1001973462'
$ws.Range("D44").Value2 = 'YAISA GONZALEZ MASCO'
$ws.Range("E44").Value2 = '2508'
$ws.Range("F44").Value2 = 56940
$ws.Range("G44").Value2 = 1423500
